# The source revision this fixture captures is a pure re-packaging /
# library-upgrade commit ("Fixed POI packaging and upgraded to POI 3.15"):
# every hunk in the change only reorders XML attributes/namespace
# declarations that Apache POI's XMLBeans serializer emits in a
# different (alphabetical) order after the upgrade. No paragraph text,
# run formatting, image, style value, or document structure actually
# changes anywhere in the diff - every "-"/"+" line pair carries the
# exact same set of attribute name/value pairs, just reshuffled.
#
# There is therefore no Word object-model mutation to make here: the
# content, styles and section properties are already correct as
# authored. We simply touch the active document so the session is
# well-formed; no edits are applied.
$d = $word.ActiveDocument
